$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCoESC-power-plants")

# Rename the existing "hydrogen" entry to "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add a new row for "hydrogen combined cycle", mirroring the formula used
# by the other power-plant rows (reference to Data!$B$20)
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Formula = "=Data!`$B`$20"

# Match the resulting selection shown in the saved workbook
$ws.Range("A24:A25").Select()
